$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Price" (column D) values. These are stored as text (e.g. "97.441.76",
# "1.80", "0.0000261") and must NOT be auto-converted to numbers, so force the cell
# format to Text before writing, then restore the default style afterwards.
$priceUpdates = @{
    'D2' = '97.441.76'
    'D3' = '3.596.60'
    'D5' = '243.61'
    'D6' = '1.80'
    'D7' = '655.03'
    'D8' = '0.426'
    'D11' = '3.594.29'
    'D15' = '4.265.11'
    'D16' = '97.311.89'
    'D17' = '0.0000261'
    'D18' = '3.579.16'
    'D20' = '12.61'
    'D21' = '18.30'
    'D22' = '0.554'
    'D23' = '3.49'
    'D24' = '517.77'
    'D26' = '6.97'
    'D27' = '103.43'
    'D28' = '13.30'
    'D29' = '0.179'
    'D32' = '1.00'
    'D34' = '0.998'
    'D35' = '31.94'
    'D38' = '617.37'
    'D39' = '8.75'
    'D41' = '1.94'
    'D44' = '6.16'
    'D45' = '0.0454'
    'D46' = '0.442'
    'D48' = '23.66'
    'D49' = '8.77'
    'D51' = '32.69'
}
foreach ($ref in $priceUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$ref]
    $cell.Style = "Normal"
}

# --- Update "Coin" (B), "Link" (C) and "Volume(1h)" (E) values. These are plain/
# percentage text and are safe to assign directly.
$textUpdates = @{
    'E2' = '  +1.86%  '
    'E3' = '  +0.47%  '
    'E4' = '  +0.00%  '
    'E5' = '  +2.51%  '
    'E6' = '  +17.15%  '
    'E7' = '  +0.14%  '
    'E8' = '  +5.95%  '
    'E9' = '  +4.73%  '
    'E10' = '  -0.04%  '
    'E11' = '  +0.44%  '
    'E12' = '  +5.06%  '
    'E13' = '  +0.52%  '
    'E14' = '  +0.34%  '
    'E15' = '  +0.49%  '
    'E16' = '  +1.86%  '
    'E17' = '  +2.24%  '
    'E18' = '  -0.28%  '
    'E19' = '  +0.29%  '
    'E20' = '  -1.21%  '
    'E21' = '  +2.04%  '
    'E22' = '  +11.00%  '
    'E23' = '  +1.21%  '
    'E24' = '  +1.20%  '
    'E25' = '  +3.09%  '
    'E26' = '  -0.74%  '
    'E27' = '  +7.97%  '
    'E28' = '  +4.04%  '
    'E29' = '  +23.86%  '
    'E30' = '  -1.35%  '
    'E31' = '  +4.55%  '
    'E32' = '  +0.01%  '
    'E33' = '  +6.36%  '
    'E34' = '  +0.02%  '
    'E35' = '  +0.05%  '
    'E36' = '  +2.70%  '
    'E37' = '  +3.74%  '
    'E38' = '  +3.73%  '
    'E39' = '  +2.45%  '
    'E40' = '  +1.56%  '
    'E41' = '  +2.89%  '
    'E43' = '  +1.77%  '
    'B44' = 'Filecoin'
    'C44' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'E44' = '  +5.96%  '
    'B45' = 'VeChain'
    'C45' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'E45' = '  +8.62%  '
    'B46' = 'Algorand'
    'C46' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'E46' = '  +36.75%  '
    'E47' = '  +1.19%  '
    'B48' = 'WhiteBITCoin'
    'C48' = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
    'E48' = '  +0.96%  '
    'B49' = 'Cosmos'
    'C49' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'E49' = '  +6.69%  '
    'E50' = '  +6.66%  '
    'E51' = '  -5.48%  '
}
foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}
